$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "OAuthSecurity" column header
$ws.Range("G1").Value = "OAuthSecurity"

# Populate OAuthSecurity column values
$ws.Range("G2").Value = "yes"
$ws.Range("G3").Value = "no"
$ws.Range("G5").Value = "no"
$ws.Range("G6").Value = "no"

# Update existing APIKeySecurity value for the google row
$ws.Range("E2").Value = "no"

# Match the saved selection state
$ws.Range("E2").Select()
